$d = $word.ActiveDocument

# -------------------------------------------------------------------
# 1) Title paragraph: insert "JOAH " before "Symbolic Calculator
#    Proposal" and move the "_GoBack" bookmark to sit right after it.
# -------------------------------------------------------------------
$titleStart = $d.Range(0, 0)
$titleStart.InsertBefore("JOAH ")

# "JOAH " is 5 characters long, so the bookmark goes right after it.
$bmRange = $d.Range(5, 5)
$d.Bookmarks.Add("_GoBack", $bmRange)

# -------------------------------------------------------------------
# 2) "Cross-functionality & portability" paragraph: drop the tab and
#    split "Java" / "script" / " which enables ... systems." into
#    three runs.
# -------------------------------------------------------------------
$d.Content.Find.Execute("across ^tmultiple", $true, $false, $false, $false, $false, $true, 1, $false, "across multiple", 2) | Out-Null

$crossRange = $d.Content
$crossRange.Find.Execute("Java which")
$crossJavaEnd = $crossRange.Start + 4

$crossIns = $d.Range($crossJavaEnd, $crossJavaEnd)
$crossIns.InsertBefore("script")

$crossScriptEnd = $crossJavaEnd + 6

$crossBmA = $d.Range($crossJavaEnd, $crossJavaEnd)
$d.Bookmarks.Add("SplitCrossA", $crossBmA)
$crossBmB = $d.Range($crossScriptEnd, $crossScriptEnd)
$d.Bookmarks.Add("SplitCrossB", $crossBmB)
$d.Bookmarks("SplitCrossA").Delete()
$d.Bookmarks("SplitCrossB").Delete()

# -------------------------------------------------------------------
# 3) "Build the basic structure of a calculation engine using Java":
#    append a new "script" run.
# -------------------------------------------------------------------
$buildRange = $d.Content
$buildRange.Find.Execute("Build the basic structure of a calculation engine using Java")
$buildJavaEnd = $buildRange.End

$buildIns = $d.Range($buildJavaEnd, $buildJavaEnd)
$buildIns.InsertAfter("script")

$buildBm = $d.Range($buildJavaEnd, $buildJavaEnd)
$d.Bookmarks.Add("SplitBuild", $buildBm)
$d.Bookmarks("SplitBuild").Delete()
